# Natmi following Dr Hou advice:
# Recomputed LR-pair stats for Rgmb-Neo1 over Sending/Target clusters {ECs, FAPs, sCs}.
# The result is a full 3x3 cluster-pair grid (rows 2-10) replacing the old 3x2 grid (rows 2-7).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs", "Rgmb", "Neo1", "ECs", 2, 0.6666666666666666, 12.531022, 37.593066, 0.1764330600309793, 0.1764330600309793, 3, 1, 4.171693666666666, 12.515081, 0.1077921033402881, 0.1077921033402881, 52.27558511426066, 470.480266028346, 0.01901809063950257, 0.01901809063950258),
  @("ECs", "Rgmb", "Neo1", "FAPs", 2, 0.6666666666666666, 12.531022, 37.593066, 0.1764330600309793, 0.1764330600309793, 3, 1, 20.39394566666667, 61.181837, 0.5269577477327285, 0.5269577477327286, 255.5569818158047, 2300.012836342242, 0.09297276793951813, 0.09297276793951816),
  @("ECs", "Rgmb", "Neo1", "sCs", 2, 0.6666666666666666, 12.531022, 37.593066, 0.1764330600309793, 0.1764330600309793, 3, 1, 14.13565266666667, 42.406958, 0.3652501489269833, 0.3652501489269833, 177.1341745503587, 1594.207570953228, 0.06444220145195856, 0.06444220145195859),
  @("FAPs", "Rgmb", "Neo1", "ECs", 3, 1, 37.18315166666667, 111.549455, 0.5235277082863637, 0.5235277082863639, 3, 1, 4.171693666666666, 12.515081, 0.1077921033402881, 0.1077921033402881, 155.1167183145394, 1396.050464830855, 0.05643215283310794, 0.05643215283310795),
  @("FAPs", "Rgmb", "Neo1", "FAPs", 3, 1, 37.18315166666667, 111.549455, 0.5235277082863637, 0.5235277082863639, 3, 1, 20.39394566666667, 61.181837, 0.5269577477327285, 0.5269577477327286, 758.3111748054262, 6824.800573248835, 0.2758769820342591, 0.2758769820342593),
  @("FAPs", "Rgmb", "Neo1", "sCs", 3, 1, 37.18315166666667, 111.549455, 0.5235277082863637, 0.5235277082863639, 3, 1, 14.13565266666667, 42.406958, 0.3652501489269833, 0.3652501489269833, 525.6081170119878, 4730.47305310789, 0.1912185734189966, 0.1912185734189967),
  @("sCs", "Rgmb", "Neo1", "ECs", 3, 1, 21.310055, 63.930165, 0.3000392316826569, 0.3000392316826569, 3, 1, 4.171693666666666, 12.515081, 0.1077921033402881, 0.1077921033402881, 88.89902147981833, 800.091193318365, 0.0323418598676776, 0.0323418598676776),
  @("sCs", "Rgmb", "Neo1", "FAPs", 3, 1, 21.310055, 63.930165, 0.3000392316826569, 0.3000392316826569, 3, 1, 20.39394566666667, 61.181837, 0.5269577477327285, 0.5269577477327286, 434.5961038236784, 3911.364934413105, 0.1581079977589512, 0.1581079977589512),
  @("sCs", "Rgmb", "Neo1", "sCs", 3, 1, 21.310055, 63.930165, 0.3000392316826569, 0.3000392316826569, 3, 1, 14.13565266666667, 42.406958, 0.3652501489269833, 0.3652501489269833, 301.2315357875634, 2711.083822088071, 0.1095893740560281, 0.1095893740560281)
)

$r = 2
foreach ($row in $rows) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c++
  }
  $r++
}